# Weekly fruit/vegetable price update: insert one new daily price record
# for "Albahaca" (basil) dated 2022-06-02, shifting all subsequent rows
# down by one (row 364 onward -> 365 onward), growing the used range
# from A1:R388 to A1:R389.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 364-388 down one position, leaving a blank row 364.
$ws.Rows(364).Insert()

# Populate the newly-inserted row 364 with the new price observation.
$ws.Range("A364").Value = 9
$ws.Range("B364").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C364").Value = "Metropolitana"
$ws.Range("D364").Value = 44714
$ws.Range("E364").Value = 13
$ws.Range("F364").Value = 100112052
$ws.Range("G364").Value = "Albahaca"
$ws.Range("H364").Value = "Sin especificar"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 520
$ws.Range("K364").Value = 4000
$ws.Range("L364").Value = 4000
$ws.Range("M364").Value = 4000
$ws.Range("N364").Value = "$/paquete"
$ws.Range("O364").Value = "Región de Arica y Parinacota"
$ws.Range("P364").Value = 4000
$ws.Range("Q364").Value = 1
$ws.Range("R364").Value = "Hortaliza"
